$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.369.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.789.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.77%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "340.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3936"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.77%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3460"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.25"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.196"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07484"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.516"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.786.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.98%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.140"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.49%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001094"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06687"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.43%  "

$ws.Range("E19").Value = "  -2.78%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9992"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.543"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.375.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.411"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.32%  "

$ws.Range("E26").Value = "  -4.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.500"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.456"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.986.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "135.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.032"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.013"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08844"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.80%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.626"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02423"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06524"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.425"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.98%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6820"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2213"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.95%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.251"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.352"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9990"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6381"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.868"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.25%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.134"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.59%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07172"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.171"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.50%  "
